$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Remove the trailing empty paragraph from the "Reverse a Array 2 ways" row's
# second cell (it currently has two paragraphs; it should only have one).
$rowCount = $t.Rows.Count
$reverseRow = $t.Rows.Item($rowCount - 1)
$linkCell = $reverseRow.Cells.Item(2)
$linkCell.Range.Paragraphs.Item($linkCell.Range.Paragraphs.Count).Range.Delete()

# Remove the entire last row ("Convert Primitive type array to List<>").
$t.Rows.Item($t.Rows.Count).Delete()
